$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J), copying the formatting from the
# existing header cell H1 ("IP") so the new header cells share the same style.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-82: add values for new columns I (I0) and J (IF).
# Each entry is @(row, I-value, J-value)
$data = @(
    @(2, 9, 9),
    @(3, 9, 9),
    @(4, 9, 9),
    @(5, 7, 7),
    @(6, 9, 9),
    @(7, 7, 7),
    @(8, 9, 9),
    @(9, 7, 7),
    @(10, 9, 9),
    @(11, 9, 9),
    @(12, 8, 9),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 9, 9),
    @(19, 9, 9),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 9, 9),
    @(24, 9, 9),
    @(25, 9, 9),
    @(26, 9, 9),
    @(27, 10, 10),
    @(28, 8, 9),
    @(29, 8, 8),
    @(30, 9, 9),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 9, 9),
    @(35, 8, 8),
    @(36, 9, 9),
    @(37, 9, 10),
    @(38, 9, 9),
    @(39, 6, 6),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 9, 9),
    @(44, 9, 9),
    @(45, 9, 9),
    @(46, 8, 8),
    @(47, 9, 9),
    @(48, 7, 7),
    @(49, 7, 7),
    @(50, 9, 9),
    @(51, 8, 9),
    @(52, 8, 8),
    @(53, 8, 8),
    @(54, 9, 9),
    @(55, 9, 9),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 9, 9),
    @(59, 9, 9),
    @(60, 9, 9),
    @(61, 9, 9),
    @(62, 10, 10),
    @(63, 9, 9),
    @(64, 9, 9),
    @(65, 9, 9),
    @(66, 9, 9),
    @(67, 9, 9),
    @(68, 9, 9),
    @(69, 9, 9),
    @(70, 9, 9),
    @(71, 8, 8),
    @(72, 9, 9),
    @(73, 7, 7),
    @(74, 7, 7),
    @(75, 5, 5),
    @(76, 9, 9),
    @(77, 7, 7),
    @(78, 6, 6),
    @(79, 5, 5),
    @(80, 8, 8),
    @(81, 5, 5),
    @(82, 6, 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
